$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Set faces on a die"
$ws.Range("B22").Value = "I want to be able to create imbalanced dice, ie 1,1,1,2,2,3"
$ws.Range("D22").Value = "Alpha PiPi - Store Review"

$ws.Range("B24").Select()
